{"js": "// Insert a rotated \"goode\" text box and an accent-colored rectangle shape\n// into the second (empty) paragraph of the document body, matching the\n// reference OOXML diff exactly. Word's Office.js API has no direct\n// \"add floating shape\" call, so we splice the authored run-level OOXML\n// (two <w:r> runs, each an mc:AlternateContent drawing) into that\n// paragraph via Range.insertOoxml (Flat OPC), which preserves the\n// paragraph's own <w:pPr> and appends the runs before the paragraph end.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document has an introductory paragraph followed by one empty\n// paragraph right before the section break \u2014 that empty paragraph is\n// where the shapes were dropped in.\nconst targetParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst targetRange = targetParagraph.getRange();\n\nconst flatOpcXml = \"<?xml version=\\\"1.0\\\" standalone=\\\"yes\\\"?>\\n<?mso-application progid=\\\"Word.Document\\\"?>\\n<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\">\\n  <pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\\n    <pkg:xmlData>\\n      <w:document xmlns:mc=\\\"http://schemas.openxmlformats.org/markup-compatibility/2006\\\" xmlns:o=\\\"urn:schemas-microsoft-com:office:office\\\" xmlns:r=\\\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\\\" xmlns:v=\\\"urn:schemas-microsoft-com:vml\\\" xmlns:wp14=\\\"http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing\\\" xmlns:wp=\\\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\\\" xmlns:w10=\\\"urn:schemas-microsoft-com:office:word\\\" xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\" xmlns:w14=\\\"http://schemas.microsoft.com/office/word/2010/wordml\\\" xmlns:wps=\\\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\\\" xmlns:a=\\\"http://schemas.openxmlformats.org/drawingml/2006/main\\\" mc:Ignorable=\\\"w14 wp14\\\">\\n        <w:body>\\n          <w:p>\\n<w:r>\\n        <w:rPr>\\n          <w:noProof/>\\n          <w:lang w:val=\\\"en-US\\\"/>\\n        </w:rPr>\\n        <mc:AlternateContent>\\n          <mc:Choice Requires=\\\"wps\\\">\\n            <w:drawing>\\n              <wp:anchor distT=\\\"0\\\" distB=\\\"0\\\" distL=\\\"114300\\\" distR=\\\"114300\\\" simplePos=\\\"0\\\" relativeHeight=\\\"251659264\\\" behindDoc=\\\"0\\\" locked=\\\"0\\\" layoutInCell=\\\"1\\\" allowOverlap=\\\"1\\\" wp14:anchorId=\\\"477713D7\\\" wp14:editId=\\\"28FB218F\\\">\\n                <wp:simplePos x=\\\"0\\\" y=\\\"0\\\"/>\\n                <wp:positionH relativeFrom=\\\"column\\\">\\n                  <wp:posOffset>824898</wp:posOffset>\\n                </wp:positionH>\\n                <wp:positionV relativeFrom=\\\"paragraph\\\">\\n                  <wp:posOffset>845802</wp:posOffset>\\n                </wp:positionV>\\n                <wp:extent cx=\\\"1958869\\\" cy=\\\"1790700\\\"/>\\n                <wp:effectExtent l=\\\"292100\\\" t=\\\"368300\\\" r=\\\"302260\\\" b=\\\"368300\\\"/>\\n                <wp:wrapNone/>\\n                <wp:docPr id=\\\"1\\\" name=\\\"Text Box 1\\\"/>\\n                <wp:cNvGraphicFramePr/>\\n                <a:graphic xmlns:a=\\\"http://schemas.openxmlformats.org/drawingml/2006/main\\\">\\n                  <a:graphicData uri=\\\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\\\">\\n                    <wps:wsp>\\n                      <wps:cNvSpPr txBox=\\\"1\\\"/>\\n                      <wps:spPr>\\n                        <a:xfrm rot=\\\"1688297\\\">\\n                          <a:off x=\\\"0\\\" y=\\\"0\\\"/>\\n                          <a:ext cx=\\\"1958869\\\" cy=\\\"1790700\\\"/>\\n                        </a:xfrm>\\n                        <a:prstGeom prst=\\\"rect\\\">\\n                          <a:avLst/>\\n                        </a:prstGeom>\\n                        <a:solidFill>\\n                          <a:schemeClr val=\\\"lt1\\\"/>\\n                        </a:solidFill>\\n                        <a:ln w=\\\"6350\\\">\\n                          <a:solidFill>\\n                            <a:prstClr val=\\\"black\\\"/>\\n                          </a:solidFill>\\n                        </a:ln>\\n                      </wps:spPr>\\n                      <wps:txbx>\\n                        <w:txbxContent>\\n                          <w:p>\\n                            <w:pPr>\\n                              <w:rPr>\\n                                <w:lang w:val=\\\"en-US\\\"/>\\n                              </w:rPr>\\n                            </w:pPr>\\n                            <w:proofErr w:type=\\\"spellStart\\\"/>\\n                            <w:r>\\n                              <w:rPr>\\n                                <w:lang w:val=\\\"en-US\\\"/>\\n                              </w:rPr>\\n                              <w:t>goode</w:t>\\n                            </w:r>\\n                            <w:proofErr w:type=\\\"spellEnd\\\"/>\\n                          </w:p>\\n                        </w:txbxContent>\\n                      </wps:txbx>\\n                      <wps:bodyPr rot=\\\"0\\\" spcFirstLastPara=\\\"0\\\" vertOverflow=\\\"overflow\\\" horzOverflow=\\\"overflow\\\" vert=\\\"horz\\\" wrap=\\\"square\\\" lIns=\\\"91440\\\" tIns=\\\"45720\\\" rIns=\\\"91440\\\" bIns=\\\"45720\\\" numCol=\\\"1\\\" spcCol=\\\"0\\\" rtlCol=\\\"0\\\" fromWordArt=\\\"0\\\" anchor=\\\"t\\\" anchorCtr=\\\"0\\\" forceAA=\\\"0\\\" compatLnSpc=\\\"1\\\">\\n                        <a:prstTxWarp prst=\\\"textNoShape\\\">\\n                          <a:avLst/>\\n                        </a:prstTxWarp>\\n                        <a:noAutofit/>\\n                      </wps:bodyPr>\\n                    </wps:wsp>\\n                  </a:graphicData>\\n                </a:graphic>\\n                <wp14:sizeRelH relativeFrom=\\\"margin\\\">\\n                  <wp14:pctWidth>0</wp14:pctWidth>\\n                </wp14:sizeRelH>\\n              </wp:anchor>\\n            </w:drawing>\\n          </mc:Choice>\\n          <mc:Fallback>\\n            <w:pict>\\n              <v:shapetype w14:anchorId=\\\"477713D7\\\" id=\\\"_x0000_t202\\\" coordsize=\\\"21600,21600\\\" o:spt=\\\"202\\\" path=\\\"m,l,21600r21600,l21600,xe\\\">\\n                <v:stroke joinstyle=\\\"miter\\\"/>\\n                <v:path gradientshapeok=\\\"t\\\" o:connecttype=\\\"rect\\\"/>\\n              </v:shapetype>\\n              <v:shape id=\\\"Text Box 1\\\" o:spid=\\\"_x0000_s1026\\\" type=\\\"#_x0000_t202\\\" style=\\\"position:absolute;margin-left:64.95pt;margin-top:66.6pt;width:154.25pt;height:141pt;rotation:1844071fd;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-width-relative:margin;v-text-anchor:top\\\" o:gfxdata=\\\"UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#13;&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#13;&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#13;&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#13;&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#13;&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#13;&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#13;&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#13;&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#13;&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#13;&#10;IQCvT8L0VQIAALAEAAAOAAAAZHJzL2Uyb0RvYy54bWysVE1v2zAMvQ/YfxB0X+1kbb4Qp8hadBhQ&#13;&#10;tAXSoWdFlhtjsqhJSuzu1+9JcdK022nYRaDI5yfykfT8sms02ynnazIFH5zlnCkjqazNc8G/P958&#13;&#10;mnDmgzCl0GRUwV+U55eLjx/mrZ2pIW1Il8oxkBg/a23BNyHYWZZ5uVGN8GdklUGwIteIgKt7zkon&#13;&#10;WrA3Ohvm+ShryZXWkVTew3u9D/JF4q8qJcN9VXkVmC44cgvpdOlcxzNbzMXs2Qm7qWWfhviHLBpR&#13;&#10;Gzx6pLoWQbCtq/+gamrpyFMVziQ1GVVVLVWqAdUM8nfVrDbCqlQLxPH2KJP/f7TybvfgWF2id5wZ&#13;&#10;0aBFj6oL7At1bBDVaa2fAbSygIUO7ojs/R7OWHRXuYY5griD0WQynI6TFCiOAQ3VX45KR2oZKaYX&#13;&#10;k8loyplEbDCe5uM89SLbk0VS63z4qqhh0Si4QysTrdjd+oAEAD1AItyTrsubWut0ieOjrrRjO4HG&#13;&#10;65BSxhdvUNqwtuCjzxd5In4Ti9TH79dayB+x6LcMuGkDZ5RoL0W0Qrfuen3WVL5AtqQMZPBW3tTg&#13;&#10;vRU+PAiHOYMTuxPucVSakAz1Fmcbcr/+5o94tB9RzlrMbcH9z61wijP9zWAwpoPz8zjo6XJ+MR7i&#13;&#10;4k4j69OI2TZXBIXQfGSXzIgP+mBWjponrNgyvoqQMBJvFzwczKuw3yasqFTLZQJhtK0It2ZlZaQ+&#13;&#10;dPOxexLO9v0MGIU7Oky4mL1r6x4bvzS03Aaq6tTzKPBe1V53rEVqS7/Cce9O7wn1+qNZ/AYAAP//&#13;&#10;AwBQSwMEFAAGAAgAAAAhACP7IoPhAAAAEAEAAA8AAABkcnMvZG93bnJldi54bWxMT8tOwzAQvCPx&#13;&#10;D9YicaPOo62aNE5VgXIE0QL3bWzigB8hdtrw9ywnuKxmtLOzM9Vutoad1Rh67wSkiwSYcq2XvesE&#13;&#10;vL40dxtgIaKTaLxTAr5VgF19fVVhKf3FHdT5GDtGJi6UKEDHOJSch1Yri2HhB+Vo9+5Hi5Ho2HE5&#13;&#10;4oXMreFZkqy5xd7RB42Duteq/TxOVkDxmPb7j2ZaveFzqhuDX082XwtxezM/bGnst8CimuPfBfx2&#13;&#10;oPxQU7CTn5wMzBDPioKkBPI8A0aKZb5ZAjsRSFcZ8Lri/4vUPwAAAP//AwBQSwECLQAUAAYACAAA&#13;&#10;ACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQIt&#13;&#10;ABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQIt&#13;&#10;ABQABgAIAAAAIQCvT8L0VQIAALAEAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBL&#13;&#10;AQItABQABgAIAAAAIQAj+yKD4QAAABABAAAPAAAAAAAAAAAAAAAAAK8EAABkcnMvZG93bnJldi54&#13;&#10;bWxQSwUGAAAAAAQABADzAAAAvQUAAAAA&#13;&#10;\\\" fillcolor=\\\"white [3201]\\\" strokeweight=\\\".5pt\\\">\\n                <v:textbox>\\n                  <w:txbxContent>\\n                    <w:p>\\n                      <w:pPr>\\n                        <w:rPr>\\n                          <w:lang w:val=\\\"en-US\\\"/>\\n                        </w:rPr>\\n                      </w:pPr>\\n                      <w:proofErr w:type=\\\"spellStart\\\"/>\\n                      <w:r>\\n                        <w:rPr>\\n                          <w:lang w:val=\\\"en-US\\\"/>\\n                        </w:rPr>\\n                        <w:t>goode</w:t>\\n                      </w:r>\\n                      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n                    </w:p>\\n                  </w:txbxContent>\\n                </v:textbox>\\n              </v:shape>\\n            </w:pict>\\n          </mc:Fallback>\\n        </mc:AlternateContent>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:noProof/>\\n          <w:lang w:val=\\\"en-US\\\"/>\\n        </w:rPr>\\n        <mc:AlternateContent>\\n          <mc:Choice Requires=\\\"wps\\\">\\n            <w:drawing>\\n              <wp:anchor distT=\\\"0\\\" distB=\\\"0\\\" distL=\\\"114300\\\" distR=\\\"114300\\\" simplePos=\\\"0\\\" relativeHeight=\\\"251660288\\\" behindDoc=\\\"0\\\" locked=\\\"0\\\" layoutInCell=\\\"1\\\" allowOverlap=\\\"1\\\" wp14:anchorId=\\\"227920BF\\\" wp14:editId=\\\"4F9F9C0A\\\">\\n                <wp:simplePos x=\\\"0\\\" y=\\\"0\\\"/>\\n                <wp:positionH relativeFrom=\\\"column\\\">\\n                  <wp:posOffset>3221990</wp:posOffset>\\n                </wp:positionH>\\n                <wp:positionV relativeFrom=\\\"paragraph\\\">\\n                  <wp:posOffset>428614</wp:posOffset>\\n                </wp:positionV>\\n                <wp:extent cx=\\\"1341620\\\" cy=\\\"1349114\\\"/>\\n                <wp:effectExtent l=\\\"0\\\" t=\\\"0\\\" r=\\\"17780\\\" b=\\\"10160\\\"/>\\n                <wp:wrapNone/>\\n                <wp:docPr id=\\\"2\\\" name=\\\"Rectangle 2\\\"/>\\n                <wp:cNvGraphicFramePr/>\\n                <a:graphic xmlns:a=\\\"http://schemas.openxmlformats.org/drawingml/2006/main\\\">\\n                  <a:graphicData uri=\\\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\\\">\\n                    <wps:wsp>\\n                      <wps:cNvSpPr/>\\n                      <wps:spPr>\\n                        <a:xfrm>\\n                          <a:off x=\\\"0\\\" y=\\\"0\\\"/>\\n                          <a:ext cx=\\\"1341620\\\" cy=\\\"1349114\\\"/>\\n                        </a:xfrm>\\n                        <a:prstGeom prst=\\\"rect\\\">\\n                          <a:avLst/>\\n                        </a:prstGeom>\\n                      </wps:spPr>\\n                      <wps:style>\\n                        <a:lnRef idx=\\\"2\\\">\\n                          <a:schemeClr val=\\\"accent1\\\">\\n                            <a:shade val=\\\"50000\\\"/>\\n                          </a:schemeClr>\\n                        </a:lnRef>\\n                        <a:fillRef idx=\\\"1\\\">\\n                          <a:schemeClr val=\\\"accent1\\\"/>\\n                        </a:fillRef>\\n                        <a:effectRef idx=\\\"0\\\">\\n                          <a:schemeClr val=\\\"accent1\\\"/>\\n                        </a:effectRef>\\n                        <a:fontRef idx=\\\"minor\\\">\\n                          <a:schemeClr val=\\\"lt1\\\"/>\\n                        </a:fontRef>\\n                      </wps:style>\\n                      <wps:bodyPr rot=\\\"0\\\" spcFirstLastPara=\\\"0\\\" vertOverflow=\\\"overflow\\\" horzOverflow=\\\"overflow\\\" vert=\\\"horz\\\" wrap=\\\"square\\\" lIns=\\\"91440\\\" tIns=\\\"45720\\\" rIns=\\\"91440\\\" bIns=\\\"45720\\\" numCol=\\\"1\\\" spcCol=\\\"0\\\" rtlCol=\\\"0\\\" fromWordArt=\\\"0\\\" anchor=\\\"ctr\\\" anchorCtr=\\\"0\\\" forceAA=\\\"0\\\" compatLnSpc=\\\"1\\\">\\n                        <a:prstTxWarp prst=\\\"textNoShape\\\">\\n                          <a:avLst/>\\n                        </a:prstTxWarp>\\n                        <a:noAutofit/>\\n                      </wps:bodyPr>\\n                    </wps:wsp>\\n                  </a:graphicData>\\n                </a:graphic>\\n              </wp:anchor>\\n            </w:drawing>\\n          </mc:Choice>\\n          <mc:Fallback>\\n            <w:pict>\\n              <v:rect w14:anchorId=\\\"3607AD94\\\" id=\\\"Rectangle 2\\\" o:spid=\\\"_x0000_s1026\\\" style=\\\"position:absolute;margin-left:253.7pt;margin-top:33.75pt;width:105.65pt;height:106.25pt;z-index:251660288;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle\\\" o:gfxdata=\\\"UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#13;&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#13;&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#13;&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#13;&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#13;&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#13;&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#13;&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#13;&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#13;&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#13;&#10;IQAzMhk+cgIAADoFAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5faTpCoOKFFUgpkkI&#13;&#10;KmDi2Th2E8n2eWe3affrd3bSgADtYVofXDt3993d5+98frGzhm0VhhZcxcujCWfKSahbt674z8fr&#13;&#10;L6echShcLQw4VfG9Cvxi8fnTeefnagoNmFohIxAX5p2veBOjnxdFkI2yIhyBV46MGtCKSEdcFzWK&#13;&#10;jtCtKaaTyUnRAdYeQaoQ6OtVb+SLjK+1kvFO66AiMxWn2mJeMa/PaS0W52K+RuGbVg5liH+oworW&#13;&#10;UdIR6kpEwTbYvoOyrUQIoOORBFuA1q1UuQfqppy86eahEV7lXoic4Eeawv+DlbfbFbK2rviUMycs&#13;&#10;XdE9kSbc2ig2TfR0PszJ68GvcDgF2qZedxpt+qcu2C5Tuh8pVbvIJH0sv87KkykxL8lGh7OynCXU&#13;&#10;4iXcY4jfFViWNhVHSp+pFNubEHvXgwvFpXL6AvIu7o1KNRh3rzT1QSmnOTorSF0aZFtBdy+kVC6W&#13;&#10;vakRteo/H0/oN9QzRuTqMmBC1q0xI/YAkNT5HruvdfBPoSoLcAye/K2wPniMyJnBxTHYtg7wIwBD&#13;&#10;XQ2Ze/8DST01iaVnqPd0ywi9/IOX1y1xfSNCXAkkvdP90AzHO1q0ga7iMOw4awB/f/Q9+ZMMycpZ&#13;&#10;R/NT8fBrI1BxZn44EuhZOZulgcuH2fG3pAF8bXl+bXEbewl0TSW9Fl7mbfKP5rDVCPaJRn2ZspJJ&#13;&#10;OEm5Ky4jHg6XsZ9reiykWi6zGw2ZF/HGPXiZwBOrSUuPuyeBfhBcJK3ewmHWxPyN7nrfFOlguYmg&#13;&#10;2yzKF14HvmlAs3CGxyS9AK/P2evlyVv8AQAA//8DAFBLAwQUAAYACAAAACEAle6jVOEAAAAPAQAA&#13;&#10;DwAAAGRycy9kb3ducmV2LnhtbExPy07DMBC8I/EP1iJxo3YqWkdpnAoVcUHi0JYPcOMlDo3tKHaa&#13;&#10;5O9ZTnAZaTWz8yj3s+vYDYfYBq8gWwlg6OtgWt8o+Dy/PeXAYtLe6C54VLBShH11f1fqwoTJH/F2&#13;&#10;Sg0jEx8LrcCm1Becx9qi03EVevTEfYXB6UTn0HAz6InMXcfXQmy5062nBKt7PFisr6fRUYjG45LJ&#13;&#10;6XD9sPN7i93yjeOi1OPD/LojeNkBSzinvw/43UD9oaJilzB6E1mnYCPkM0kVbOUGGAlklktgFwXr&#13;&#10;XAjgVcn/76h+AAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAA&#13;&#10;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAA&#13;&#10;AAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhADMyGT5yAgAAOgUAAA4AAAAA&#13;&#10;AAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAJXuo1ThAAAADwEAAA8A&#13;&#10;AAAAAAAAAAAAAAAAzAQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAADaBQAAAAA=&#13;&#10;\\\" fillcolor=\\\"#4472c4 [3204]\\\" strokecolor=\\\"#1f3763 [1604]\\\" strokeweight=\\\"1pt\\\"/>\\n            </w:pict>\\n          </mc:Fallback>\\n        </mc:AlternateContent>\\n      </w:r>\\n          </w:p>\\n        </w:body>\\n      </w:document>\\n    </pkg:xmlData>\\n  </pkg:part>\\n</pkg:package>\";\n\ntargetRange.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Insert a rotated \"goode\" text box and an accent-colored rectangle shape\n# into the second (empty) paragraph of the document body, matching the\n# reference OOXML diff exactly. The Word COM object model has no direct\n# \"add floating shape from scratch\" call that reproduces this authored\n# markup, so we splice the authored run-level OOXML (two <w:r> runs, each\n# an mc:AlternateContent drawing) into that paragraph's Range via\n# Range.InsertXML (Flat OPC payload), targeting the *end* of the existing\n# paragraph range so the paragraph's own mark/pPr survive untouched and\n# the runs land right before it, exactly like the captured edit.\n\n$d = $word.ActiveDocument\n\n# The document has an introductory paragraph followed by one empty\n# paragraph right before the section break -- that empty paragraph is\n# where the shapes were dropped in.\n$targetParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$targetRange = $targetParagraph.Range\n\n$flatOpcXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:mc=\"http://schemas.openxmlformats.org/markup-compatibility/2006\" xmlns:o=\"urn:schemas-microsoft-com:office:office\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" xmlns:v=\"urn:schemas-microsoft-com:vml\" xmlns:wp14=\"http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing\" xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" xmlns:w10=\"urn:schemas-microsoft-com:office:word\" xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" xmlns:wps=\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\" xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" mc:Ignorable=\"w14 wp14\">\n        <w:body>\n          <w:p>\n<w:r>\n        <w:rPr>\n          <w:noProof/>\n          <w:lang w:val=\"en-US\"/>\n        </w:rPr>\n        <mc:AlternateContent>\n          <mc:Choice Requires=\"wps\">\n            <w:drawing>\n              <wp:anchor distT=\"0\" distB=\"0\" distL=\"114300\" distR=\"114300\" simplePos=\"0\" relativeHeight=\"251659264\" behindDoc=\"0\" locked=\"0\" layoutInCell=\"1\" allowOverlap=\"1\" wp14:anchorId=\"477713D7\" wp14:editId=\"28FB218F\">\n                <wp:simplePos x=\"0\" y=\"0\"/>\n                <wp:positionH relativeFrom=\"column\">\n                  <wp:posOffset>824898</wp:posOffset>\n                </wp:positionH>\n                <wp:positionV relativeFrom=\"paragraph\">\n                  <wp:posOffset>845802</wp:posOffset>\n                </wp:positionV>\n                <wp:extent cx=\"1958869\" cy=\"1790700\"/>\n                <wp:effectExtent l=\"292100\" t=\"368300\" r=\"302260\" b=\"368300\"/>\n                <wp:wrapNone/>\n                <wp:docPr id=\"1\" name=\"Text Box 1\"/>\n                <wp:cNvGraphicFramePr/>\n                <a:graphic xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\">\n                  <a:graphicData uri=\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\">\n                    <wps:wsp>\n                      <wps:cNvSpPr txBox=\"1\"/>\n                      <wps:spPr>\n                        <a:xfrm rot=\"1688297\">\n                          <a:off x=\"0\" y=\"0\"/>\n                          <a:ext cx=\"1958869\" cy=\"1790700\"/>\n                        </a:xfrm>\n                        <a:prstGeom prst=\"rect\">\n                          <a:avLst/>\n                        </a:prstGeom>\n                        <a:solidFill>\n                          <a:schemeClr val=\"lt1\"/>\n                        </a:solidFill>\n                        <a:ln w=\"6350\">\n                          <a:solidFill>\n                            <a:prstClr val=\"black\"/>\n                          </a:solidFill>\n                        </a:ln>\n                      </wps:spPr>\n                      <wps:txbx>\n                        <w:txbxContent>\n                          <w:p>\n                            <w:pPr>\n                              <w:rPr>\n                                <w:lang w:val=\"en-US\"/>\n                              </w:rPr>\n                            </w:pPr>\n                            <w:proofErr w:type=\"spellStart\"/>\n                            <w:r>\n                              <w:rPr>\n                                <w:lang w:val=\"en-US\"/>\n                              </w:rPr>\n                              <w:t>goode</w:t>\n                            </w:r>\n                            <w:proofErr w:type=\"spellEnd\"/>\n                          </w:p>\n                        </w:txbxContent>\n                      </wps:txbx>\n                      <wps:bodyPr rot=\"0\" spcFirstLastPara=\"0\" vertOverflow=\"overflow\" horzOverflow=\"overflow\" vert=\"horz\" wrap=\"square\" lIns=\"91440\" tIns=\"45720\" rIns=\"91440\" bIns=\"45720\" numCol=\"1\" spcCol=\"0\" rtlCol=\"0\" fromWordArt=\"0\" anchor=\"t\" anchorCtr=\"0\" forceAA=\"0\" compatLnSpc=\"1\">\n                        <a:prstTxWarp prst=\"textNoShape\">\n                          <a:avLst/>\n                        </a:prstTxWarp>\n                        <a:noAutofit/>\n                      </wps:bodyPr>\n                    </wps:wsp>\n                  </a:graphicData>\n                </a:graphic>\n                <wp14:sizeRelH relativeFrom=\"margin\">\n                  <wp14:pctWidth>0</wp14:pctWidth>\n                </wp14:sizeRelH>\n              </wp:anchor>\n            </w:drawing>\n          </mc:Choice>\n          <mc:Fallback>\n            <w:pict>\n              <v:shapetype w14:anchorId=\"477713D7\" id=\"_x0000_t202\" coordsize=\"21600,21600\" o:spt=\"202\" path=\"m,l,21600r21600,l21600,xe\">\n                <v:stroke joinstyle=\"miter\"/>\n                <v:path gradientshapeok=\"t\" o:connecttype=\"rect\"/>\n              </v:shapetype>\n              <v:shape id=\"Text Box 1\" o:spid=\"_x0000_s1026\" type=\"#_x0000_t202\" style=\"position:absolute;margin-left:64.95pt;margin-top:66.6pt;width:154.25pt;height:141pt;rotation:1844071fd;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-width-relative:margin;v-text-anchor:top\" o:gfxdata=\"UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#13;&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#13;&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#13;&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#13;&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#13;&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#13;&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#13;&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#13;&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#13;&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#13;&#10;IQCvT8L0VQIAALAEAAAOAAAAZHJzL2Uyb0RvYy54bWysVE1v2zAMvQ/YfxB0X+1kbb4Qp8hadBhQ&#13;&#10;tAXSoWdFlhtjsqhJSuzu1+9JcdK022nYRaDI5yfykfT8sms02ynnazIFH5zlnCkjqazNc8G/P958&#13;&#10;mnDmgzCl0GRUwV+U55eLjx/mrZ2pIW1Il8oxkBg/a23BNyHYWZZ5uVGN8GdklUGwIteIgKt7zkon&#13;&#10;WrA3Ohvm+ShryZXWkVTew3u9D/JF4q8qJcN9VXkVmC44cgvpdOlcxzNbzMXs2Qm7qWWfhviHLBpR&#13;&#10;Gzx6pLoWQbCtq/+gamrpyFMVziQ1GVVVLVWqAdUM8nfVrDbCqlQLxPH2KJP/f7TybvfgWF2id5wZ&#13;&#10;0aBFj6oL7At1bBDVaa2fAbSygIUO7ojs/R7OWHRXuYY5griD0WQynI6TFCiOAQ3VX45KR2oZKaYX&#13;&#10;k8loyplEbDCe5uM89SLbk0VS63z4qqhh0Si4QysTrdjd+oAEAD1AItyTrsubWut0ieOjrrRjO4HG&#13;&#10;65BSxhdvUNqwtuCjzxd5In4Ti9TH79dayB+x6LcMuGkDZ5RoL0W0Qrfuen3WVL5AtqQMZPBW3tTg&#13;&#10;vRU+PAiHOYMTuxPucVSakAz1Fmcbcr/+5o94tB9RzlrMbcH9z61wijP9zWAwpoPz8zjo6XJ+MR7i&#13;&#10;4k4j69OI2TZXBIXQfGSXzIgP+mBWjponrNgyvoqQMBJvFzwczKuw3yasqFTLZQJhtK0It2ZlZaQ+&#13;&#10;dPOxexLO9v0MGIU7Oky4mL1r6x4bvzS03Aaq6tTzKPBe1V53rEVqS7/Cce9O7wn1+qNZ/AYAAP//&#13;&#10;AwBQSwMEFAAGAAgAAAAhACP7IoPhAAAAEAEAAA8AAABkcnMvZG93bnJldi54bWxMT8tOwzAQvCPx&#13;&#10;D9YicaPOo62aNE5VgXIE0QL3bWzigB8hdtrw9ywnuKxmtLOzM9Vutoad1Rh67wSkiwSYcq2XvesE&#13;&#10;vL40dxtgIaKTaLxTAr5VgF19fVVhKf3FHdT5GDtGJi6UKEDHOJSch1Yri2HhB+Vo9+5Hi5Ho2HE5&#13;&#10;4oXMreFZkqy5xd7RB42Duteq/TxOVkDxmPb7j2ZaveFzqhuDX082XwtxezM/bGnst8CimuPfBfx2&#13;&#10;oPxQU7CTn5wMzBDPioKkBPI8A0aKZb5ZAjsRSFcZ8Lri/4vUPwAAAP//AwBQSwECLQAUAAYACAAA&#13;&#10;ACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQIt&#13;&#10;ABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQIt&#13;&#10;ABQABgAIAAAAIQCvT8L0VQIAALAEAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBL&#13;&#10;AQItABQABgAIAAAAIQAj+yKD4QAAABABAAAPAAAAAAAAAAAAAAAAAK8EAABkcnMvZG93bnJldi54&#13;&#10;bWxQSwUGAAAAAAQABADzAAAAvQUAAAAA&#13;&#10;\" fillcolor=\"white [3201]\" strokeweight=\".5pt\">\n                <v:textbox>\n                  <w:txbxContent>\n                    <w:p>\n                      <w:pPr>\n                        <w:rPr>\n                          <w:lang w:val=\"en-US\"/>\n                        </w:rPr>\n                      </w:pPr>\n                      <w:proofErr w:type=\"spellStart\"/>\n                      <w:r>\n                        <w:rPr>\n                          <w:lang w:val=\"en-US\"/>\n                        </w:rPr>\n                        <w:t>goode</w:t>\n                      </w:r>\n                      <w:proofErr w:type=\"spellEnd\"/>\n                    </w:p>\n                  </w:txbxContent>\n                </v:textbox>\n              </v:shape>\n            </w:pict>\n          </mc:Fallback>\n        </mc:AlternateContent>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:noProof/>\n          <w:lang w:val=\"en-US\"/>\n        </w:rPr>\n        <mc:AlternateContent>\n          <mc:Choice Requires=\"wps\">\n            <w:drawing>\n              <wp:anchor distT=\"0\" distB=\"0\" distL=\"114300\" distR=\"114300\" simplePos=\"0\" relativeHeight=\"251660288\" behindDoc=\"0\" locked=\"0\" layoutInCell=\"1\" allowOverlap=\"1\" wp14:anchorId=\"227920BF\" wp14:editId=\"4F9F9C0A\">\n                <wp:simplePos x=\"0\" y=\"0\"/>\n                <wp:positionH relativeFrom=\"column\">\n                  <wp:posOffset>3221990</wp:posOffset>\n                </wp:positionH>\n                <wp:positionV relativeFrom=\"paragraph\">\n                  <wp:posOffset>428614</wp:posOffset>\n                </wp:positionV>\n                <wp:extent cx=\"1341620\" cy=\"1349114\"/>\n                <wp:effectExtent l=\"0\" t=\"0\" r=\"17780\" b=\"10160\"/>\n                <wp:wrapNone/>\n                <wp:docPr id=\"2\" name=\"Rectangle 2\"/>\n                <wp:cNvGraphicFramePr/>\n                <a:graphic xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\">\n                  <a:graphicData uri=\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\">\n                    <wps:wsp>\n                      <wps:cNvSpPr/>\n                      <wps:spPr>\n                        <a:xfrm>\n                          <a:off x=\"0\" y=\"0\"/>\n                          <a:ext cx=\"1341620\" cy=\"1349114\"/>\n                        </a:xfrm>\n                        <a:prstGeom prst=\"rect\">\n                          <a:avLst/>\n                        </a:prstGeom>\n                      </wps:spPr>\n                      <wps:style>\n                        <a:lnRef idx=\"2\">\n                          <a:schemeClr val=\"accent1\">\n                            <a:shade val=\"50000\"/>\n                          </a:schemeClr>\n                        </a:lnRef>\n                        <a:fillRef idx=\"1\">\n                          <a:schemeClr val=\"accent1\"/>\n                        </a:fillRef>\n                        <a:effectRef idx=\"0\">\n                          <a:schemeClr val=\"accent1\"/>\n                        </a:effectRef>\n                        <a:fontRef idx=\"minor\">\n                          <a:schemeClr val=\"lt1\"/>\n                        </a:fontRef>\n                      </wps:style>\n                      <wps:bodyPr rot=\"0\" spcFirstLastPara=\"0\" vertOverflow=\"overflow\" horzOverflow=\"overflow\" vert=\"horz\" wrap=\"square\" lIns=\"91440\" tIns=\"45720\" rIns=\"91440\" bIns=\"45720\" numCol=\"1\" spcCol=\"0\" rtlCol=\"0\" fromWordArt=\"0\" anchor=\"ctr\" anchorCtr=\"0\" forceAA=\"0\" compatLnSpc=\"1\">\n                        <a:prstTxWarp prst=\"textNoShape\">\n                          <a:avLst/>\n                        </a:prstTxWarp>\n                        <a:noAutofit/>\n                      </wps:bodyPr>\n                    </wps:wsp>\n                  </a:graphicData>\n                </a:graphic>\n              </wp:anchor>\n            </w:drawing>\n          </mc:Choice>\n          <mc:Fallback>\n            <w:pict>\n              <v:rect w14:anchorId=\"3607AD94\" id=\"Rectangle 2\" o:spid=\"_x0000_s1026\" style=\"position:absolute;margin-left:253.7pt;margin-top:33.75pt;width:105.65pt;height:106.25pt;z-index:251660288;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle\" o:gfxdata=\"UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#13;&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#13;&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#13;&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#13;&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#13;&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#13;&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#13;&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#13;&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#13;&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#13;&#10;IQAzMhk+cgIAADoFAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5faTpCoOKFFUgpkkI&#13;&#10;KmDi2Th2E8n2eWe3affrd3bSgADtYVofXDt3993d5+98frGzhm0VhhZcxcujCWfKSahbt674z8fr&#13;&#10;L6echShcLQw4VfG9Cvxi8fnTeefnagoNmFohIxAX5p2veBOjnxdFkI2yIhyBV46MGtCKSEdcFzWK&#13;&#10;jtCtKaaTyUnRAdYeQaoQ6OtVb+SLjK+1kvFO66AiMxWn2mJeMa/PaS0W52K+RuGbVg5liH+oworW&#13;&#10;UdIR6kpEwTbYvoOyrUQIoOORBFuA1q1UuQfqppy86eahEV7lXoic4Eeawv+DlbfbFbK2rviUMycs&#13;&#10;XdE9kSbc2ig2TfR0PszJ68GvcDgF2qZedxpt+qcu2C5Tuh8pVbvIJH0sv87KkykxL8lGh7OynCXU&#13;&#10;4iXcY4jfFViWNhVHSp+pFNubEHvXgwvFpXL6AvIu7o1KNRh3rzT1QSmnOTorSF0aZFtBdy+kVC6W&#13;&#10;vakRteo/H0/oN9QzRuTqMmBC1q0xI/YAkNT5HruvdfBPoSoLcAye/K2wPniMyJnBxTHYtg7wIwBD&#13;&#10;XQ2Ze/8DST01iaVnqPd0ywi9/IOX1y1xfSNCXAkkvdP90AzHO1q0ga7iMOw4awB/f/Q9+ZMMycpZ&#13;&#10;R/NT8fBrI1BxZn44EuhZOZulgcuH2fG3pAF8bXl+bXEbewl0TSW9Fl7mbfKP5rDVCPaJRn2ZspJJ&#13;&#10;OEm5Ky4jHg6XsZ9reiykWi6zGw2ZF/HGPXiZwBOrSUuPuyeBfhBcJK3ewmHWxPyN7nrfFOlguYmg&#13;&#10;2yzKF14HvmlAs3CGxyS9AK/P2evlyVv8AQAA//8DAFBLAwQUAAYACAAAACEAle6jVOEAAAAPAQAA&#13;&#10;DwAAAGRycy9kb3ducmV2LnhtbExPy07DMBC8I/EP1iJxo3YqWkdpnAoVcUHi0JYPcOMlDo3tKHaa&#13;&#10;5O9ZTnAZaTWz8yj3s+vYDYfYBq8gWwlg6OtgWt8o+Dy/PeXAYtLe6C54VLBShH11f1fqwoTJH/F2&#13;&#10;Sg0jEx8LrcCm1Becx9qi03EVevTEfYXB6UTn0HAz6InMXcfXQmy5062nBKt7PFisr6fRUYjG45LJ&#13;&#10;6XD9sPN7i93yjeOi1OPD/LojeNkBSzinvw/43UD9oaJilzB6E1mnYCPkM0kVbOUGGAlklktgFwXr&#13;&#10;XAjgVcn/76h+AAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAA&#13;&#10;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAA&#13;&#10;AAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhADMyGT5yAgAAOgUAAA4AAAAA&#13;&#10;AAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAJXuo1ThAAAADwEAAA8A&#13;&#10;AAAAAAAAAAAAAAAAzAQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAADaBQAAAAA=&#13;&#10;\" fillcolor=\"#4472c4 [3204]\" strokecolor=\"#1f3763 [1604]\" strokeweight=\"1pt\"/>\n            </w:pict>\n          </mc:Fallback>\n        </mc:AlternateContent>\n      </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$targetRange.InsertXML($flatOpcXml, \"End\")\n"}
